# Community - aggregator control remaining capacity of households
#
# Inserts a new column "H" on Sheet1 holding the new
# "aggregator_household_battery_control" flag (value 1 for every scenario
# row), shifting the existing price/id columns one slot to the right
# (old H:K -> I:L), and updates the sheet's zoom / selection to match the
# author's final view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H; this shifts the old H:K columns (buy/sell
# price factor, id_electricity, id_electricity_feed_in) to I:L and carries
# their column widths/formatting along automatically.
$ws.Columns("H").Insert()

# Match the width of the neighbouring price-factor columns (F:G).
$ws.Columns("H").ColumnWidth = 23

# Header for the new column.
$ws.Range("H1").Value = "aggregator_household_battery_control"

# Every scenario row gets the flag enabled.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Restore the view state captured in the saved workbook.
[void]$ws.Range("H7").Select()
$excel.ActiveWindow.Zoom = 137
